$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 9).Value = 1.0
$ws.Cells.Item(5, 9).Value = 1.0
$ws.Cells.Item(6, 11).Value = "Yes"
$ws.Cells.Item(7, 7).Value = 2.0
$ws.Cells.Item(7, 9).Value = 1.0
$ws.Cells.Item(8, 9).Value = 1.0
$ws.Cells.Item(8, 11).Value = "Yes"
$ws.Cells.Item(9, 9).Value = 1.0
$ws.Cells.Item(9, 11).Value = "Yes"
$ws.Cells.Item(11, 9).Value = 1.0
$ws.Cells.Item(11, 11).Value = "Yes"
$ws.Cells.Item(13, 9).Value = 1.0
$ws.Cells.Item(14, 9).Value = 1.0
$ws.Cells.Item(16, 7).Value = 1.0
$ws.Cells.Item(19, 9).Value = 1.0
$ws.Cells.Item(24, 9).Value = 1.0
$ws.Cells.Item(24, 11).Value = "Yes"
$ws.Cells.Item(25, 11).Value = "Yes"
$ws.Cells.Item(28, 11).Value = "Yes"
$ws.Cells.Item(30, 9).Value = 1.0
$ws.Cells.Item(33, 11).Value = "Yes"
$ws.Cells.Item(36, 7).Value = 2.0
$ws.Cells.Item(36, 9).Value = 1.0
$ws.Cells.Item(39, 9).Value = 1.0
$ws.Cells.Item(43, 9).Value = 1.0
$ws.Cells.Item(43, 11).Value = "Yes"
$ws.Cells.Item(47, 9).Value = 1.0
$ws.Cells.Item(49, 9).Value = 1.0
$ws.Cells.Item(52, 9).Value = 1.0
$ws.Cells.Item(53, 9).Value = 1.0
$ws.Cells.Item(54, 7).Value = 1.0
$ws.Cells.Item(54, 9).Value = 1.0
$ws.Cells.Item(54, 11).Value = "Yes"
$ws.Cells.Item(55, 9).Value = 1.0
$ws.Cells.Item(55, 11).Value = "Yes"
$ws.Cells.Item(58, 9).Value = 1.0
$ws.Cells.Item(58, 11).Value = "Yes"
$ws.Cells.Item(59, 9).Value = 1.0
$ws.Cells.Item(59, 11).Value = "Yes"
$ws.Cells.Item(61, 9).Value = 1.0
$ws.Cells.Item(62, 9).Value = 1.0
$ws.Cells.Item(64, 9).Value = 1.0
$ws.Cells.Item(64, 11).Value = "Yes"
$ws.Cells.Item(66, 9).Value = 1.0
$ws.Cells.Item(69, 9).Value = 1.0
$ws.Cells.Item(69, 11).Value = "Yes"
$ws.Cells.Item(72, 9).Value = 1.0
$ws.Cells.Item(73, 9).Value = 1.0
$ws.Cells.Item(73, 11).Value = "Yes"
$ws.Cells.Item(74, 9).Value = 1.0
$ws.Cells.Item(75, 11).Value = "Yes"
$ws.Cells.Item(76, 9).Value = 1.0
$ws.Cells.Item(77, 9).Value = 1.0
$ws.Cells.Item(78, 7).Value = 1.0
$ws.Cells.Item(78, 9).Value = 1.0
$ws.Cells.Item(80, 9).Value = 1.0
$ws.Cells.Item(86, 9).Value = 1.0
$ws.Cells.Item(86, 11).Value = "Yes"
$ws.Cells.Item(87, 11).Value = "Yes"
$ws.Cells.Item(88, 9).Value = 1.0
$ws.Cells.Item(89, 9).Value = 1.0
$ws.Cells.Item(93, 11).Value = "Yes"
$ws.Cells.Item(94, 9).Value = 1.0
$ws.Cells.Item(95, 9).Value = 1.0
$ws.Cells.Item(96, 9).Value = 1.0
$ws.Cells.Item(96, 11).Value = "Yes"
$ws.Cells.Item(100, 9).Value = 1.0
$ws.Cells.Item(103, 7).Value = 3.0
$ws.Cells.Item(110, 7).Value = 1.0
$ws.Cells.Item(110, 9).Value = 1.0
$ws.Cells.Item(110, 11).Value = "Yes"
$ws.Cells.Item(111, 7).Value = 4.0
$ws.Cells.Item(111, 8).Value = 3.0
$ws.Cells.Item(111, 9).Value = 1.0
$ws.Cells.Item(111, 11).Value = "Yes"
$ws.Cells.Item(112, 9).Value = 1.0
$ws.Cells.Item(112, 11).Value = "Yes"
$ws.Cells.Item(121, 9).Value = 1.0
$ws.Cells.Item(121, 11).Value = "Yes"
$ws.Cells.Item(122, 9).Value = 1.0
$ws.Cells.Item(124, 9).Value = 1.0
$ws.Cells.Item(124, 11).Value = "Yes"
$ws.Cells.Item(125, 9).Value = 1.0
$ws.Cells.Item(125, 11).Value = "Yes"
$ws.Cells.Item(126, 9).Value = 1.0
$ws.Cells.Item(127, 9).Value = 1.0
$ws.Cells.Item(127, 11).Value = "Yes"
$ws.Cells.Item(136, 9).Value = 1.0
$ws.Cells.Item(137, 9).Value = 1.0
$ws.Cells.Item(143, 9).Value = 1.0
$ws.Cells.Item(145, 7).Value = 2.0
$ws.Cells.Item(150, 9).Value = 1.0
$ws.Cells.Item(154, 9).Value = 1.0
$ws.Cells.Item(157, 9).Value = 1.0
$ws.Cells.Item(158, 9).Value = 1.0
$ws.Cells.Item(160, 9).Value = 1.0
$ws.Cells.Item(161, 9).Value = 1.0
$ws.Cells.Item(163, 11).Value = "Yes"
$ws.Cells.Item(169, 11).Value = "Yes"
$ws.Cells.Item(170, 9).Value = 1.0
$ws.Cells.Item(173, 9).Value = 1.0
$ws.Cells.Item(174, 9).Value = 1.0
$ws.Cells.Item(177, 9).Value = 1.0
